$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while preserving the underlying
# Text cell type (matches how the source workbook stores all data cells,
# including Price, as inline/shared strings rather than numbers).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "246.31"

# Row 3
Set-TextValue $ws.Range("D3") "22.05"

# Row 4
Set-TextValue $ws.Range("D4") "5.447"

# Row 5
Set-TextValue $ws.Range("D5") "0.05779"

# Row 7
Set-TextValue $ws.Range("D7") "6.317"

# Row 8
Set-TextValue $ws.Range("D8") "0.8183"

# Row 9
Set-TextValue $ws.Range("D9") "0.9699"
$ws.Range("E9").Value = "8FTXTokenFTT"

# Row 11
Set-TextValue $ws.Range("D11") "0.07506"

# Row 12
Set-TextValue $ws.Range("D12") "0.03126"

# Row 13
Set-TextValue $ws.Range("D13") "0.02999"

# Row 14
Set-TextValue $ws.Range("D14") "4.151"

# Row 15
Set-TextValue $ws.Range("D15") "0.09405"

# Row 16
Set-TextValue $ws.Range("D16") "0.001589"

# Row 17
Set-TextValue $ws.Range("D17") "0.04804"

# Row 19
Set-TextValue $ws.Range("D19") "0.006193"

# Row 21
Set-TextValue $ws.Range("D21") "0.0009958"

# Row 23
Set-TextValue $ws.Range("D23") "3.766"

# Row 25
Set-TextValue $ws.Range("D25") "0.3229"

# Row 27
Set-TextValue $ws.Range("D27") "0.0003999"

# Row 40
Set-TextValue $ws.Range("D40") "0.03893"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006461"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1075"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003000"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

# Row 44
Set-TextValue $ws.Range("D44") "0.006632"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005592"

# Row 47
Set-TextValue $ws.Range("D47") "0.3800"
